$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 9,20
$data[0,0] = "ECs"
$data[0,1] = "Fgf16"
$data[0,2] = "Fgfr3"
$data[0,3] = "ECs"
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.09476766666666665
$data[0,7] = 0.284303
$data[0,8] = 0.02456723071903196
$data[0,9] = 0.02456723071903196
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 2.619953333333334
$data[0,13] = 7.85986
$data[0,14] = 0.6100029538328192
$data[0,15] = 0.6100029538328192
$data[0,16] = 0.2482868641755555
$data[0,17] = 2.23458177758
$data[0,18] = 0.01498608330610187
$data[0,19] = 0.01498608330610187

$data[1,0] = "ECs"
$data[1,1] = "Fgf16"
$data[1,2] = "Fgfr3"
$data[1,3] = "FAPs"
$data[1,4] = 1
$data[1,5] = 0.3333333333333333
$data[1,6] = 0.09476766666666665
$data[1,7] = 0.284303
$data[1,8] = 0.02456723071903196
$data[1,9] = 0.02456723071903196
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.5698483333333333
$data[1,13] = 1.709545
$data[1,14] = 0.1326776176306101
$data[1,15] = 0.1326776176306101
$data[1,16] = 0.05400319690388888
$data[1,17] = 0.4860287721349999
$data[1,18] = 0.0032595216435827
$data[1,19] = 0.0032595216435827

$data[2,0] = "ECs"
$data[2,1] = "Fgf16"
$data[2,2] = "Fgfr3"
$data[2,3] = "sCs"
$data[2,4] = 1
$data[2,5] = 0.3333333333333333
$data[2,6] = 0.09476766666666665
$data[2,7] = 0.284303
$data[2,8] = 0.02456723071903196
$data[2,9] = 0.02456723071903196
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 1.105183
$data[2,13] = 3.315549
$data[2,14] = 0.2573194285365706
$data[2,15] = 0.2573194285365706
$data[2,16] = 0.1047356141496667
$data[2,17] = 0.9426205273469999
$data[2,18] = 0.006321625769347387
$data[2,19] = 0.006321625769347388

$data[3,0] = "FAPs"
$data[3,1] = "Fgf16"
$data[3,2] = "Fgfr3"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 2.187396333333333
$data[3,7] = 6.562189
$data[3,8] = 0.567052796435119
$data[3,9] = 0.5670527964351191
$data[3,10] = 2
$data[3,11] = 0.6666666666666666
$data[3,12] = 2.619953333333334
$data[3,13] = 7.85986
$data[3,14] = 0.6100029538328192
$data[3,15] = 0.6100029538328192
$data[3,16] = 5.730876314837778
$data[3,17] = 51.57788683354
$data[3,18] = 0.3459038808045829
$data[3,19] = 0.345903880804583

$data[4,0] = "FAPs"
$data[4,1] = "Fgf16"
$data[4,2] = "Fgfr3"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 2.187396333333333
$data[4,7] = 6.562189
$data[4,8] = 0.567052796435119
$data[4,9] = 0.5670527964351191
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.5698483333333333
$data[4,13] = 1.709545
$data[4,14] = 0.1326776176306101
$data[4,15] = 0.1326776176306101
$data[4,16] = 1.246484154889445
$data[4,17] = 11.218357394005
$data[4,18] = 0.07523521410178688
$data[4,19] = 0.07523521410178689

$data[5,0] = "FAPs"
$data[5,1] = "Fgf16"
$data[5,2] = "Fgfr3"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 2.187396333333333
$data[5,7] = 6.562189
$data[5,8] = 0.567052796435119
$data[5,9] = 0.5670527964351191
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 1.105183
$data[5,13] = 3.315549
$data[5,14] = 0.2573194285365706
$data[5,15] = 0.2573194285365706
$data[5,16] = 2.417473241862333
$data[5,17] = 21.757259176761
$data[5,18] = 0.1459137015287491
$data[5,19] = 0.1459137015287492

$data[6,0] = "sCs"
$data[6,1] = "Fgf16"
$data[6,2] = "Fgfr3"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 1.575318666666667
$data[6,7] = 4.725956
$data[6,8] = 0.408379972845849
$data[6,9] = 0.408379972845849
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 2.619953333333334
$data[6,13] = 7.85986
$data[6,14] = 0.6100029538328192
$data[6,15] = 0.6100029538328192
$data[6,16] = 4.127261391795556
$data[6,17] = 37.14535252616
$data[6,18] = 0.2491129897221344
$data[6,19] = 0.2491129897221344

$data[7,0] = "sCs"
$data[7,1] = "Fgf16"
$data[7,2] = "Fgfr3"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 1.575318666666667
$data[7,7] = 4.725956
$data[7,8] = 0.408379972845849
$data[7,9] = 0.408379972845849
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 0.5698483333333333
$data[7,13] = 1.709545
$data[7,14] = 0.1326776176306101
$data[7,15] = 0.1326776176306101
$data[7,16] = 0.8976927166688888
$data[7,17] = 8.07923445002
$data[7,18] = 0.05418288188524048
$data[7,19] = 0.05418288188524048

$data[8,0] = "sCs"
$data[8,1] = "Fgf16"
$data[8,2] = "Fgfr3"
$data[8,3] = "sCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 1.575318666666667
$data[8,7] = 4.725956
$data[8,8] = 0.408379972845849
$data[8,9] = 0.408379972845849
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 1.105183
$data[8,13] = 3.315549
$data[8,14] = 0.2573194285365706
$data[8,15] = 0.2573194285365706
$data[8,16] = 1.741015409982667
$data[8,17] = 15.669138689844
$data[8,18] = 0.1050841012384741
$data[8,19] = 0.1050841012384741

$ws.Range("A2:T10").Value = $data
